$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Half Span): rows 2-5 become formulas doubling the original value
$ws.Range("A2").Formula = "=40*2"
$ws.Range("A3").Formula = "=37*2"
$ws.Range("A4").Formula = "=35*2"
$ws.Range("A5").Formula = "=33*2"

# Column E (Deflection): rows 4-7 get their integer part replaced with 4,
# keeping the original fractional part
$e4old = 1.296
$e5old = 1.052
$e6old = 0.75509999999999999
$e7old = 0.022179999999999998

$ws.Range("E4").Value = 4 + ($e4old - [math]::Floor($e4old))
$ws.Range("E5").Value = 4 + ($e5old - [math]::Floor($e5old))
$ws.Range("E6").Value = 4 + ($e6old - [math]::Floor($e6old))
$ws.Range("E7").Value = 4 + ($e7old - [math]::Floor($e7old))

# Update the selected cell on the sheet
$ws.Range("F12").Select()
